$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new hours value for D23 (this feeds the SUM formula in F3)
$ws.Range("D23").Value = 8

# Force recalculation so the cached formula value updates (86.5 -> 94.5)
$excel.Calculate()

# Update the view/selection to match target state: no frozen topLeftCell,
# selection on D23 instead of F20
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("D23").Select()
